$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$line1 = "O que vai ser preciso para criar as mensagens nulas"
$line2 = "Utilizar a média de envio das mensagens da aplicação para calcular a promessa"
$line3 = "Criar pseudoAlgoritmo do envio das mensagens nulas e da média das mensagens. "
$line4 = "Verificar a possibilidade de realizar o cálculo do checkpoint sem a utilização da mensagem nula"
$line5 = "DCB tem que criar tudo, aplicação apenas faz operações internas. "

# Build five bullet paragraphs in one shot - this inherits the bullet
# pPr / black solidFill rPr that already lived on the placeholder's
# single (blank) run.
$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5

# Paragraph 3 needs to be split into three runs: "Criar ", "pseudoAlgoritmo",
# " do envio das mensagens nulas e da média das mensagens. "
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "Criar "
$run3b = $para3.InsertAfter("pseudoAlgoritmo")
$run3c = $run3b.InsertAfter(" do envio das mensagens nulas e da média das mensagens. ")

# Paragraph 5 needs to be split into two runs: "DCB tem que criar tudo,
# aplicação apenas faz " and "operações internas. "
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "DCB tem que criar tudo, aplicação apenas faz "
$run5b = $para5.InsertAfter("operações internas. ")
